# Update detected GPS coordinates (column B) for the "api_detected_locations"
# sheet as part of the "final spider web routes fallback" change: several
# rows' coordinates were corrected / replaced with fallback values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "32.767884,34.966961"
$ws.Range("B8").Value = "33.084756,35.112133"
$ws.Range("B11").Value = "32.128872,34.948863"
$ws.Range("B18").Value = "31.854920,35.218710"
$ws.Range("B19").Value = "32.019045,34.841227"
$ws.Range("B21").Value = "32.267628,34.993511"
$ws.Range("B24").Value = "32.773251,35.044543"
$ws.Range("B25").Value = "32.606459,35.290914"
$ws.Range("B29").Value = "31.854920,35.218710"
$ws.Range("B30").Value = "31.961622,34.807607"
$ws.Range("B31").Value = "32.035964,34.845985"
$ws.Range("B32").Value = "32.267628,34.993511"
$ws.Range("B34").Value = "31.767207,35.224441"
$ws.Range("B35").Value = "32.052244,34.797643"
$ws.Range("B36").Value = "32.047035,34.899314"
$ws.Range("B37").Value = "31.790191,35.198620"
$ws.Range("B41").Value = "32.233023,34.950420"
$ws.Range("B42").Value = "32.267628,34.993511"
$ws.Range("B43").Value = "32.049544,34.764454"
$ws.Range("B51").Value = "32.037040,34.776415"
$ws.Range("B56").Value = "32.081982,34.816659"
$ws.Range("B57").Value = "32.028209,34.802593"
$ws.Range("B58").Value = "33.132610,35.690627"
$ws.Range("B61").Value = "31.767207,35.224441"
$ws.Range("B64").Value = "32.064156,34.854185"
$ws.Range("B69").Value = "32.153195,34.846595"
$ws.Range("B70").Value = "32.006200,34.743653"
$ws.Range("B75").Value = "32.037040,34.776415"
$ws.Range("B77").Value = "31.790191,35.198620"
$ws.Range("B83").Value = "32.028209,34.802593"
$ws.Range("B85").Value = "31.944800,34.877389"
$ws.Range("B88").Value = "32.099281,34.896845"
$ws.Range("B91").Value = "31.755751,34.983774"
$ws.Range("B93").Value = "31.750988,35.207798"
$ws.Range("B95").Value = "31.853707,35.217433"
$ws.Range("B96").Value = "31.194371,34.837706"
$ws.Range("B97").Value = "31.750585,35.215673"
$ws.Range("B100").Value = "31.223027,34.809387"
$ws.Range("B101").Value = "31.065689,35.014440"
$ws.Range("B102").Value = "31.241723,34.804322"
$ws.Range("B104").Value = "32.174844,34.814576"
$ws.Range("B105").Value = "31.864522,34.741564"
$ws.Range("B106").Value = "31.927646,34.878243"
$ws.Range("B107").Value = "32.175016,34.928954"
$ws.Range("B108").Value = "32.045844,34.752383"
$ws.Range("B109").Value = "31.785936,35.221741"
$ws.Range("B111").Value = "31.785077,34.693905"
$ws.Range("B112").Value = "32.472006,34.946602"
$ws.Range("B113").Value = "31.667321,34.601532"
$ws.Range("B114").Value = "32.267563,34.993779"
$ws.Range("B116").Value = "31.975998,34.882170"
$ws.Range("B119").Value = "31.669726,34.779153"
$ws.Range("B120").Value = "32.092353,34.885480"
$ws.Range("B121").Value = "31.244467,34.807280"
$ws.Range("B123").Value = "31.757029,34.990864"
$ws.Range("B124").Value = "32.053835,34.771023"

Write-Host "Updated 56 coordinate cells"
